# Snake leaderboard: append the new game results recorded after adding
# walls to the game (Sheet1 holds the combined/overall leaderboard).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Name, Score pairs to append starting at row 6 (rows 1-5 already exist).
$entries = @(
    @("Jack", 1042.0),
    @("test", 2084.0),
    @("jg",   1042.0),
    @("jg",   3647.0),
    @("a",    521.0),
    @("jg",   521.0),
    @("jh",   521.0),
    @("jk",   4689.0),
    @("jk",   4689.0)
)

$row = 6
foreach ($entry in $entries) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
